# Ineos Known Locomotive List - remove "Parent company" and
# "Location County/City" columns (additional companies sent for
# questionaire; these columns are no longer needed on this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("Location County/City") first so column B's index
# is not affected by the later deletion.
$ws.Columns.Item(5).Delete()

# Delete column B ("Parent company").
$ws.Columns.Item(2).Delete()

$ws.Range("M5").Select()
